# Delete row 456 ("「毎年、ノーベル賞が授けられる」" post) from Sheet1.
# This shifts all subsequent rows up by one, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(456).EntireRow.Delete()
